$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new shared string / header for column F
$ws.Range("F1").Value = "FRAZ_PLURICAND"

# Fill F2:F15 with 0.2 (20%)
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 6).Value = 0.2
}

# Apply the same number format / style as column E to the new column F
$ws.Range("F1:F15").NumberFormat = $ws.Range("E1:E15").NumberFormat
$ws.Columns.Item(6).ColumnWidth = 15.83

# Update selection to match the diff (active cell F14)
$ws.Range("F14").Select()
